$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on C29 (it pointed to http://ideone.com/g23yCk),
# without disturbing the other sheet hyperlinks (B2, A10).
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq "`$C`$29") {
        $h.Delete()
    }
}

# C29 now gets a new URL (plain text, no hyperlink)
$ws.Range("C29").Value = "http://ideone.com/gv4dZS"

# New row 30: "Delete Tree" problem, with the old ideone link moved down (plain text)
$ws.Range("A30").Value = "Delete Tree"
$ws.Range("C30").Value = "http://ideone.com/g23yCk"

# Match column C styling (Hyperlink-style formatting used throughout col B/C)
$ws.Range("C30").Style = "Hyperlink"

# Update the selection to match post-edit state
$ws.Range("B31").Select()
